$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.27"
$ws.Range("G2").Value = "'7"
$ws.Range("D3").Value = "'29.35"
$ws.Range("E3").Value = "'-1.90%"
$ws.Range("G3").Value = "'7"
$ws.Range("D4").Value = "'5.137"
$ws.Range("E4").Value = "'-0.34%"
$ws.Range("G4").Value = "'7"
$ws.Range("D5").Value = "'0.05799"
$ws.Range("E5").Value = "'2.21%"
$ws.Range("G5").Value = "'7"
$ws.Range("D6").Value = "'6.660"
$ws.Range("E6").Value = "'1.75%"
$ws.Range("G6").Value = "'7"
$ws.Range("D7").Value = "'3.233"
$ws.Range("E7").Value = "'6.77%"
$ws.Range("G7").Value = "'7"
$ws.Range("D8").Value = "'0.8520"
$ws.Range("E8").Value = "'0.30%"
$ws.Range("G8").Value = "'7"
$ws.Range("D9").Value = "'0.8599"
$ws.Range("E9").Value = "'-0.92%"
$ws.Range("G9").Value = "'7"
$ws.Range("D10").Value = "'0.1382"
$ws.Range("E10").Value = "'2.83%"
$ws.Range("G10").Value = "'7"
$ws.Range("D11").Value = "'0.07092"
$ws.Range("E11").Value = "'2.62%"
$ws.Range("G11").Value = "'7"
$ws.Range("D12").Value = "'0.03170"
$ws.Range("E12").Value = "'9.70%"
$ws.Range("G12").Value = "'7"
$ws.Range("D13").Value = "'0.09378"
$ws.Range("E13").Value = "'0.02%"
$ws.Range("G13").Value = "'7"
$ws.Range("D14").Value = "'0.001531"
$ws.Range("E14").Value = "'1.14%"
$ws.Range("G14").Value = "'7"
$ws.Range("D15").Value = "'0.0005978"
$ws.Range("E15").Value = "'-94.06%"
$ws.Range("G15").Value = "'7"
$ws.Range("D16").Value = "'0.006150"
$ws.Range("E16").Value = "'3.04%"
$ws.Range("G16").Value = "'7"
$ws.Range("D17").Value = "'3.494"
$ws.Range("E17").Value = "'-0.41%"
$ws.Range("G17").Value = "'7"
$ws.Range("D18").Value = "'2.224"
$ws.Range("E18").Value = "'-0.82%"
$ws.Range("G18").Value = "'7"
$ws.Range("D19").Value = "'0.3196"
$ws.Range("E19").Value = "'1.51%"
$ws.Range("G19").Value = "'7"
$ws.Range("D20").Value = "'0.03352"
$ws.Range("E20").Value = "'0.05%"
$ws.Range("G20").Value = "'7"
$ws.Range("E21").Value = "'-1.61%"
$ws.Range("G21").Value = "'7"
$ws.Range("D22").Value = "'3.477"
$ws.Range("E22").Value = "'-3.65%"
$ws.Range("G22").Value = "'7"
$ws.Range("D23").Value = "'0.04135"
$ws.Range("E23").Value = "'-1.14%"
$ws.Range("G23").Value = "'7"
$ws.Range("D24").Value = "'0.1381"
$ws.Range("E24").Value = "'0.51%"
$ws.Range("G24").Value = "'7"
$ws.Range("D25").Value = "'0.001227"
$ws.Range("E25").Value = "'1.47%"
$ws.Range("G25").Value = "'7"
$ws.Range("D26").Value = "'0.004142"
$ws.Range("E26").Value = "'-6.87%"
$ws.Range("G26").Value = "'7"
$ws.Range("D27").Value = "'0.0001099"
$ws.Range("E27").Value = "'-6.78%"
$ws.Range("G27").Value = "'7"
$ws.Range("D28").Value = "'0.0001445"
$ws.Range("E28").Value = "'3.90%"
$ws.Range("G28").Value = "'7"
$ws.Range("G29").Value = "'7"
$ws.Range("G30").Value = "'7"
$ws.Range("G31").Value = "'7"
$ws.Range("G32").Value = "'7"
$ws.Range("G33").Value = "'7"
$ws.Range("G34").Value = "'7"
$ws.Range("G35").Value = "'7"
$ws.Range("G36").Value = "'7"
$ws.Range("G37").Value = "'7"
$ws.Range("G38").Value = "'7"
$ws.Range("G39").Value = "'7"
$ws.Range("D40").Value = "'0.03752"
$ws.Range("E40").Value = "'-1.12%"
$ws.Range("G40").Value = "'7"
$ws.Range("D41").Value = "'0.005757"
$ws.Range("E41").Value = "'-2.37%"
$ws.Range("G41").Value = "'7"
$ws.Range("D42").Value = "'0.1069"
$ws.Range("E42").Value = "'1.14%"
$ws.Range("G42").Value = "'7"
$ws.Range("E43").Value = "'-3.90%"
$ws.Range("G43").Value = "'7"
$ws.Range("D44").Value = "'0.009259"
$ws.Range("E44").Value = "'-0.89%"
$ws.Range("G44").Value = "'7"
$ws.Range("D45").Value = "'0.00005293"
$ws.Range("E45").Value = "'4.47%"
$ws.Range("G45").Value = "'7"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("G46").Value = "'7"
$ws.Range("D47").Value = "'0.05798"
$ws.Range("E47").Value = "'-35.54%"
$ws.Range("G47").Value = "'7"
$ws.Range("D48").Value = "'0.002178"
$ws.Range("E48").Value = "'-21.14%"
$ws.Range("G48").Value = "'7"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("G49").Value = "'7"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("G50").Value = "'7"
$ws.Range("G51").Value = "'7"
